$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "data as of" timestamp in A1
$ws.Range("A1").Value = "Datos actualizados a 2 de Julio de 2020 a las 23:17"

$ws.Cells.Item(4, 2).Value = 2823345
$ws.Cells.Item(4, 3).Value = 43392
$ws.Cells.Item(4, 4).Value = 1184464
$ws.Cells.Item(4, 5).Value = 1507494
$ws.Cells.Item(4, 7).Value = 589
$ws.Cells.Item(4, 8).Value = 131387

$ws.Cells.Item(7, 2).Value = 627168
$ws.Cells.Item(7, 3).Value = 21948
$ws.Cells.Item(7, 4).Value = 379902
$ws.Cells.Item(7, 5).Value = 229041
$ws.Cells.Item(7, 7).Value = 377
$ws.Cells.Item(7, 8).Value = 18225

$ws.Cells.Item(18, 2).Value = 196674
$ws.Cells.Item(18, 3).Value = 350
$ws.Cells.Item(18, 5).Value = 7311

$ws.Cells.Item(26, 2).Value = 71299
$ws.Cells.Item(26, 3).Value = 1485
$ws.Cells.Item(26, 4).Value = 19288
$ws.Cells.Item(26, 5).Value = 48891
$ws.Cells.Item(26, 7).Value = 86
$ws.Cells.Item(26, 8).Value = 3120

$ws.Cells.Item(31, 1).Value = "Ecuador"
$ws.Cells.Item(31, 2).Value = 59468
$ws.Cells.Item(31, 3).Value = 1211
$ws.Cells.Item(31, 4).Value = 28032
$ws.Cells.Item(31, 5).Value = 26797
$ws.Cells.Item(31, 7).Value = 63
$ws.Cells.Item(31, 8).Value = 4639

$ws.Cells.Item(32, 1).Value = "Indonesia"
$ws.Cells.Item(32, 2).Value = 59394
$ws.Cells.Item(32, 3).Value = 1624
$ws.Cells.Item(32, 4).Value = 26667
$ws.Cells.Item(32, 5).Value = 29740
$ws.Cells.Item(32, 7).Value = 53
$ws.Cells.Item(32, 8).Value = 2987

$ws.Cells.Item(51, 1).Value = "Israel"
$ws.Cells.Item(51, 2).Value = 27047
$ws.Cells.Item(51, 3).Value = 790
$ws.Cells.Item(51, 4).Value = 17547
$ws.Cells.Item(51, 5).Value = 9176
$ws.Cells.Item(51, 7).Value = 2
$ws.Cells.Item(51, 8).Value = 324

$ws.Cells.Item(52, 1).Value = "Armenia"
$ws.Cells.Item(52, 2).Value = 26658
$ws.Cells.Item(52, 3).Value = 593
$ws.Cells.Item(52, 4).Value = 15036
$ws.Cells.Item(52, 5).Value = 11163
$ws.Cells.Item(52, 7).Value = 6
$ws.Cells.Item(52, 8).Value = 459

$ws.Cells.Item(53, 1).Value = "Nigeria"
$ws.Cells.Item(53, 2).Value = 26484
$ws.Cells.Item(53, 3).Value = 0
$ws.Cells.Item(53, 4).Value = 10152
$ws.Cells.Item(53, 5).Value = 15729
$ws.Cells.Item(53, 7).Value = 0
$ws.Cells.Item(53, 8).Value = 603

$ws.Cells.Item(54, 2).Value = 25489
$ws.Cells.Item(54, 3).Value = 12
$ws.Cells.Item(54, 5).Value = 387

$ws.Cells.Item(70, 2).Value = 9992
$ws.Cells.Item(70, 3).Value = 290
$ws.Cells.Item(70, 4).Value = 4660
$ws.Cells.Item(70, 5).Value = 5264

$ws.Cells.Item(88, 2).Value = 5450
$ws.Cells.Item(88, 3).Value = 46
$ws.Cells.Item(88, 4).Value = 4392

$ws.Cells.Item(92, 2).Value = 4606
$ws.Cells.Item(92, 3).Value = 134
$ws.Cells.Item(92, 4).Value = 1727
$ws.Cells.Item(92, 5).Value = 2750

$ws.Cells.Item(96, 1).Value = "Costa Rica"
$ws.Cells.Item(96, 2).Value = 4023
$ws.Cells.Item(96, 3).Value = 270
$ws.Cells.Item(96, 4).Value = 1589
$ws.Cells.Item(96, 5).Value = 2417
$ws.Cells.Item(96, 8).Value = 17

$ws.Cells.Item(97, 1).Value = "Republica de Africa Central"
$ws.Cells.Item(97, 2).Value = 3788
$ws.Cells.Item(97, 3).Value = 43
$ws.Cells.Item(97, 4).Value = 810
$ws.Cells.Item(97, 5).Value = 2931
$ws.Cells.Item(97, 8).Value = 47

$ws.Cells.Item(100, 2).Value = 3080
$ws.Cells.Item(100, 3).Value = 322
$ws.Cells.Item(100, 5).Value = 2612

$ws.Cells.Item(109, 2).Value = 2303
$ws.Cells.Item(109, 3).Value = 43
$ws.Cells.Item(109, 4).Value = 1108
$ws.Cells.Item(109, 5).Value = 1176

$ws.Cells.Item(128, 1).Value = "Yemen"
$ws.Cells.Item(128, 2).Value = 1221
$ws.Cells.Item(128, 3).Value = 31
$ws.Cells.Item(128, 4).Value = 513
$ws.Cells.Item(128, 5).Value = 383
$ws.Cells.Item(128, 7).Value = 7
$ws.Cells.Item(128, 8).Value = 325

$ws.Cells.Item(129, 1).Value = "Benin"
$ws.Cells.Item(129, 2).Value = 1199
$ws.Cells.Item(129, 4).Value = 333
$ws.Cells.Item(129, 5).Value = 845
$ws.Cells.Item(129, 8).Value = 21

$ws.Cells.Item(134, 2).Value = 1063
$ws.Cells.Item(134, 3).Value = 21
$ws.Cells.Item(134, 4).Value = 493
$ws.Cells.Item(134, 5).Value = 567

$ws.Cells.Item(146, 2).Value = 717
$ws.Cells.Item(146, 3).Value = 2
$ws.Cells.Item(146, 4).Value = 260
$ws.Cells.Item(146, 5).Value = 444

$ws.Cells.Item(151, 2).Value = 667
$ws.Cells.Item(151, 3).Value = 6
$ws.Cells.Item(151, 4).Value = 424
$ws.Cells.Item(151, 5).Value = 229

$ws.Cells.Item(152, 1).Value = "Zimbabue"
$ws.Cells.Item(152, 2).Value = 617
$ws.Cells.Item(152, 3).Value = 12
$ws.Cells.Item(152, 4).Value = 173
$ws.Cells.Item(152, 5).Value = 437
$ws.Cells.Item(152, 8).Value = 7

$ws.Cells.Item(153, 1).Value = "Montenegro"
$ws.Cells.Item(153, 2).Value = 616
$ws.Cells.Item(153, 3).Value = 40
$ws.Cells.Item(153, 4).Value = 315
$ws.Cells.Item(153, 5).Value = 289
$ws.Cells.Item(153, 8).Value = 12

$ws.Cells.Item(161, 1).Value = "Angola"
$ws.Cells.Item(161, 2).Value = 315
$ws.Cells.Item(161, 3).Value = 24
$ws.Cells.Item(161, 4).Value = 97
$ws.Cells.Item(161, 5).Value = 201
$ws.Cells.Item(161, 7).Value = 2
$ws.Cells.Item(161, 8).Value = 17

$ws.Cells.Item(162, 1).Value = "Siria"
$ws.Cells.Item(162, 2).Value = 312
$ws.Cells.Item(162, 3).Value = 19
$ws.Cells.Item(162, 4).Value = 113
$ws.Cells.Item(162, 5).Value = 190
$ws.Cells.Item(162, 8).Value = 9

$ws.Cells.Item(163, 1).Value = "Birmania"
$ws.Cells.Item(163, 2).Value = 304
$ws.Cells.Item(163, 3).Value = 1
$ws.Cells.Item(163, 4).Value = 223
$ws.Cells.Item(163, 5).Value = 75
$ws.Cells.Item(163, 8).Value = 6

$ws.Cells.Item(164, 1).Value = "Comoras"
$ws.Cells.Item(164, 2).Value = 303
$ws.Cells.Item(164, 3).Value = 0
$ws.Cells.Item(164, 4).Value = 200
$ws.Cells.Item(164, 5).Value = 96
$ws.Cells.Item(164, 8).Value = 7

$ws.Cells.Item(165, 1).Value = "Namibia"
$ws.Cells.Item(165, 2).Value = 293
$ws.Cells.Item(165, 3).Value = 8
$ws.Cells.Item(165, 4).Value = 24
$ws.Cells.Item(165, 5).Value = 269
$ws.Cells.Item(165, 8).Value = 0

$ws.Cells.Item(190, 2).Value = 55
$ws.Cells.Item(190, 3).Value = 6
$ws.Cells.Item(190, 5).Value = 26

$ws.Cells.Item(205, 1).Value = "Fiyi"

$ws.Cells.Item(206, 1).Value = "Dominica"
